# feat: sops Update 3
#
# Rename the main form sheet from "F-SW-FR-03" (New Feature Request form)
# to "S-SW-SC-03" (Software Service Catalog form), and keep the workbook's
# Print_Area defined name pointing at the renamed sheet.

$wb = $excel.ActiveWorkbook

# The visible/active sheet holding the form (first tab, "F-SW-FR-03").
$ws = $wb.ActiveSheet

# Rename the sheet.
$ws.Name = "S-SW-SC-03"

# Excel does not auto-update the Print_Area defined name's sheet qualifier
# on a rename, so keep it in sync with the new sheet name explicitly (the
# area itself, A1:D23, is unchanged).
$ws.PageSetup.PrintArea = '$A$1:$D$23'
